$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows with corrected / swapped values ---
# Row 29
$ws.Range("B29").Value = 6865295
$ws.Range("F29").Value = "FK Tuzla City"
$ws.Range("G29").Value = "NK Igman Konjic"
$ws.Range("H29").Value = 3
$ws.Range("I29").Value = 1
$ws.Range("K29").Value = 1.8
$ws.Range("L29").Value = 3.4
$ws.Range("M29").Value = 3.8
$ws.Range("N29").Value = 1.615
$ws.Range("O29").Value = 3.5
$ws.Range("P29").Value = 4.5
$ws.Range("Q29").Value = -0.75
$ws.Range("R29").Value = 1.85
$ws.Range("S29").Value = 1.95
$ws.Range("T29").Value = 2.75
$ws.Range("U29").Value = 2
$ws.Range("V29").Value = 1.8
$ws.Range("W29").Value = 0.615
$ws.Range("Z29").Value = 0.8500000000000001
$ws.Range("AB29").Value = 1
$ws.Range("AC29").Value = -1

# Row 30
$ws.Range("B30").Value = 6865296
$ws.Range("F30").Value = "Velez Mostar"
$ws.Range("G30").Value = "Zeljeznicar"
$ws.Range("H30").Value = 1
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 1.909
$ws.Range("L30").Value = 3.2
$ws.Range("M30").Value = 3.6
$ws.Range("N30").Value = 1.95
$ws.Range("O30").Value = 3.2
$ws.Range("P30").Value = 3.4
$ws.Range("Q30").Value = -0.5
$ws.Range("R30").Value = 2.025
$ws.Range("S30").Value = 1.775
$ws.Range("T30").Value = 2.25
$ws.Range("U30").Value = 1.9
$ws.Range("V30").Value = 1.9
$ws.Range("W30").Value = 0.95
$ws.Range("Z30").Value = 1.025
$ws.Range("AB30").Value = -1
$ws.Range("AC30").Value = 0.8999999999999999

# Row 36
$ws.Range("B36").Value = 6864629
$ws.Range("F36").Value = "Borac Banja Luka"
$ws.Range("G36").Value = "NK Posusje"
$ws.Range("H36").Value = 1
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 1.363
$ws.Range("L36").Value = 4.5
$ws.Range("M36").Value = 6.5
$ws.Range("N36").Value = 1.363
$ws.Range("O36").Value = 4.2
$ws.Range("P36").Value = 6.5
$ws.Range("R36").Value = 1.95
$ws.Range("S36").Value = 1.85
$ws.Range("T36").Value = 2.5
$ws.Range("U36").Value = 1.925
$ws.Range("V36").Value = 1.875
$ws.Range("W36").Value = 0.363
$ws.Range("AA36").Value = 0.425
$ws.Range("AB36").Value = -1
$ws.Range("AC36").Value = 0.875

# Row 37
$ws.Range("B37").Value = 6865299
$ws.Range("F37").Value = "Siroki Brijeg"
$ws.Range("G37").Value = "Zvijezda 09"
$ws.Range("H37").Value = 2
$ws.Range("I37").Value = 1
$ws.Range("K37").Value = 1.25
$ws.Range("L37").Value = 5.5
$ws.Range("M37").Value = 8
$ws.Range("N37").Value = 1.4
$ws.Range("O37").Value = 4.75
$ws.Range("P37").Value = 5.75
$ws.Range("R37").Value = 1.9
$ws.Range("S37").Value = 1.9
$ws.Range("T37").Value = 2.75
$ws.Range("U37").Value = 1.85
$ws.Range("V37").Value = 1.95
$ws.Range("W37").Value = 0.3999999999999999
$ws.Range("AA37").Value = 0.45
$ws.Range("AB37").Value = 0.425
$ws.Range("AC37").Value = -0.5

# Row 76
$ws.Range("B76").Value = 6865377
$ws.Range("F76").Value = "Zrinjski Mostar"
$ws.Range("G76").Value = "FK Tuzla City"
$ws.Range("H76").Value = 3
$ws.Range("J76").Value = "H"
$ws.Range("K76").Value = 1.333
$ws.Range("L76").Value = 5
$ws.Range("M76").Value = 6
$ws.Range("N76").Value = 1.166
$ws.Range("O76").Value = 6.5
$ws.Range("P76").Value = 13
$ws.Range("Q76").Value = -2
$ws.Range("R76").Value = 1.9
$ws.Range("S76").Value = 1.9
$ws.Range("T76").Value = 3.25
$ws.Range("U76").Value = 1.95
$ws.Range("V76").Value = 1.85
$ws.Range("W76").Value = 0.1659999999999999
$ws.Range("X76").Value = -1
$ws.Range("Z76").Value = 0
$ws.Range("AA76").Value = -0
$ws.Range("AB76").Value = 0.95
$ws.Range("AC76").Value = -1

# Row 77
$ws.Range("B77").Value = 6865328
$ws.Range("F77").Value = "Siroki Brijeg"
$ws.Range("G77").Value = "NK Posusje"
$ws.Range("H77").Value = 1
$ws.Range("J77").Value = "D"
$ws.Range("K77").Value = 2
$ws.Range("L77").Value = 3
$ws.Range("M77").Value = 3.5
$ws.Range("N77").Value = 2.1
$ws.Range("O77").Value = 3
$ws.Range("P77").Value = 3.3
$ws.Range("Q77").Value = -0.25
$ws.Range("R77").Value = 1.825
$ws.Range("S77").Value = 1.975
$ws.Range("T77").Value = 2
$ws.Range("U77").Value = 1.825
$ws.Range("V77").Value = 1.975
$ws.Range("W77").Value = -1
$ws.Range("X77").Value = 2
$ws.Range("Z77").Value = -0.5
$ws.Range("AA77").Value = 0.4875
$ws.Range("AB77").Value = 0
$ws.Range("AC77").Value = -0

# Row 122
$ws.Range("B122").Value = 6865363
$ws.Range("F122").Value = "NK Igman Konjic"
$ws.Range("G122").Value = "Siroki Brijeg"
$ws.Range("H122").Value = 1
$ws.Range("K122").Value = 2
$ws.Range("L122").Value = 3.3
$ws.Range("M122").Value = 3.25
$ws.Range("N122").Value = 2.3
$ws.Range("O122").Value = 3.2
$ws.Range("P122").Value = 2.75
$ws.Range("Q122").Value = -0.25
$ws.Range("R122").Value = 2.05
$ws.Range("S122").Value = 1.75
$ws.Range("T122").Value = 2
$ws.Range("U122").Value = 1.9
$ws.Range("V122").Value = 1.9
$ws.Range("W122").Value = 1.3
$ws.Range("Z122").Value = 1.05
$ws.Range("AC122").Value = 0.8999999999999999

# Row 123
$ws.Range("B123").Value = 6865381
$ws.Range("F123").Value = "FK Tuzla City"
$ws.Range("G123").Value = "Zvijezda 09"
$ws.Range("H123").Value = 2
$ws.Range("K123").Value = 1.666
$ws.Range("L123").Value = 3.6
$ws.Range("M123").Value = 4.333
$ws.Range("N123").Value = 1.5
$ws.Range("O123").Value = 4
$ws.Range("P123").Value = 5.25
$ws.Range("Q123").Value = -1
$ws.Range("R123").Value = 1.925
$ws.Range("S123").Value = 1.875
$ws.Range("T123").Value = 2.5
$ws.Range("U123").Value = 1.8
$ws.Range("V123").Value = 2
$ws.Range("W123").Value = 0.5
$ws.Range("Z123").Value = 0.925
$ws.Range("AC123").Value = 1

# Row 128
$ws.Range("B128").Value = 6865367
$ws.Range("E128").Value = 45360.375
$ws.Range("F128").Value = "Siroki Brijeg"
$ws.Range("G128").Value = "Sloga"
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 1
$ws.Range("J128").Value = "A"
$ws.Range("K128").Value = 1.95
$ws.Range("L128").Value = 3.2
$ws.Range("M128").Value = 3.4
$ws.Range("N128").Value = 1.65
$ws.Range("O128").Value = 3.5
$ws.Range("P128").Value = 4.5
$ws.Range("Q128").Value = -0.75
$ws.Range("R128").Value = 1.925
$ws.Range("S128").Value = 1.875
$ws.Range("U128").Value = 1.75
$ws.Range("V128").Value = 2.05
$ws.Range("W128").Value = -1
$ws.Range("X128").Value = -1
$ws.Range("Y128").Value = 3.5
$ws.Range("Z128").Value = -1
$ws.Range("AA128").Value = 0.875
$ws.Range("AB128").Value = -1
$ws.Range("AC128").Value = 1.05

# Row 129
$ws.Range("B129").Value = 6865368
$ws.Range("E129").Value = 45360.5
$ws.Range("F129").Value = "GOSK Gabela"
$ws.Range("G129").Value = "NK Posusje"
$ws.Range("H129").Value = 1
$ws.Range("I129").Value = 1
$ws.Range("J129").Value = "D"
$ws.Range("K129").Value = 3.2
$ws.Range("L129").Value = 3.1
$ws.Range("M129").Value = 2.1
$ws.Range("N129").Value = 3
$ws.Range("O129").Value = 3.1
$ws.Range("P129").Value = 2.2
$ws.Range("Q129").Value = 0.25
$ws.Range("R129").Value = 1.825
$ws.Range("S129").Value = 1.975
$ws.Range("T129").Value = 2
$ws.Range("U129").Value = 1.825
$ws.Range("V129").Value = 1.975
$ws.Range("W129").Value = -1
$ws.Range("X129").Value = 2.1
$ws.Range("Y129").Value = -1
$ws.Range("Z129").Value = 0.4125
$ws.Range("AA129").Value = -0.5
$ws.Range("AB129").Value = 0
$ws.Range("AC129").Value = -0

# --- Append brand-new rows 130-133 ---
# Row 130
$ws.Range("A130").Value = 128
$ws.Range("B130").Value = 6865366
$ws.Range("C130").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D130").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E130").Value = 45360.60416666666
$ws.Range("F130").Value = "FK Sarajevo"
$ws.Range("G130").Value = "NK Igman Konjic"
$ws.Range("H130").Value = 2
$ws.Range("I130").Value = 2
$ws.Range("J130").Value = "D"
$ws.Range("K130").Value = 1.166
$ws.Range("L130").Value = 6
$ws.Range("M130").Value = 11
$ws.Range("N130").Value = 1.363
$ws.Range("O130").Value = 4.5
$ws.Range("P130").Value = 6.5
$ws.Range("Q130").Value = -1.25
$ws.Range("R130").Value = 1.875
$ws.Range("S130").Value = 1.925
$ws.Range("T130").Value = 2.75
$ws.Range("U130").Value = 1.925
$ws.Range("V130").Value = 1.875
$ws.Range("W130").Value = -1
$ws.Range("X130").Value = 3.5
$ws.Range("Y130").Value = -1
$ws.Range("Z130").Value = -1
$ws.Range("AA130").Value = 0.925
$ws.Range("AB130").Value = 0.925
$ws.Range("AC130").Value = -1
$ws.Range("A130").Font.Bold = $true
$ws.Range("A130").HorizontalAlignment = -4108
$ws.Range("A130").VerticalAlignment = -4160
$ws.Range("A130").Borders.LineStyle = 1
$ws.Range("E130").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 131
$ws.Range("A131").Value = 129
$ws.Range("B131").Value = 6865365
$ws.Range("C131").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D131").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E131").Value = 45361.41666666666
$ws.Range("F131").Value = "Zvijezda 09"
$ws.Range("G131").Value = "Zeljeznicar"
$ws.Range("H131").Value = 2
$ws.Range("I131").Value = 1
$ws.Range("J131").Value = "H"
$ws.Range("K131").Value = 2.5
$ws.Range("L131").Value = 3
$ws.Range("M131").Value = 2.75
$ws.Range("N131").Value = 2.5
$ws.Range("O131").Value = 2.875
$ws.Range("P131").Value = 2.875
$ws.Range("Q131").Value = 0
$ws.Range("R131").Value = 1.775
$ws.Range("S131").Value = 2.025
$ws.Range("T131").Value = 1.75
$ws.Range("U131").Value = 1.75
$ws.Range("V131").Value = 2.05
$ws.Range("W131").Value = 1.5
$ws.Range("X131").Value = -1
$ws.Range("Y131").Value = -1
$ws.Range("Z131").Value = 0.7749999999999999
$ws.Range("AA131").Value = -1
$ws.Range("AB131").Value = 0.75
$ws.Range("AC131").Value = -1
$ws.Range("A131").Font.Bold = $true
$ws.Range("A131").HorizontalAlignment = -4108
$ws.Range("A131").VerticalAlignment = -4160
$ws.Range("A131").Borders.LineStyle = 1
$ws.Range("E131").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 132
$ws.Range("A132").Value = 130
$ws.Range("B132").Value = 6864644
$ws.Range("C132").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D132").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E132").Value = 45361.625
$ws.Range("F132").Value = "Zrinjski Mostar"
$ws.Range("G132").Value = "Borac Banja Luka"
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 1
$ws.Range("J132").Value = "A"
$ws.Range("K132").Value = 1.8
$ws.Range("L132").Value = 3.3
$ws.Range("M132").Value = 4
$ws.Range("N132").Value = 1.85
$ws.Range("O132").Value = 3.3
$ws.Range("P132").Value = 3.75
$ws.Range("Q132").Value = -0.5
$ws.Range("R132").Value = 1.85
$ws.Range("S132").Value = 1.95
$ws.Range("T132").Value = 2
$ws.Range("U132").Value = 1.775
$ws.Range("V132").Value = 2.025
$ws.Range("W132").Value = -1
$ws.Range("X132").Value = -1
$ws.Range("Y132").Value = 2.75
$ws.Range("Z132").Value = -1
$ws.Range("AA132").Value = 0.95
$ws.Range("AB132").Value = -1
$ws.Range("AC132").Value = 1.025
$ws.Range("A132").Font.Bold = $true
$ws.Range("A132").HorizontalAlignment = -4108
$ws.Range("A132").VerticalAlignment = -4160
$ws.Range("A132").Borders.LineStyle = 1
$ws.Range("E132").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 133
$ws.Range("A133").Value = 131
$ws.Range("B133").Value = 6865369
$ws.Range("C133").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D133").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E133").Value = 45362.58333333334
$ws.Range("F133").Value = "Velez Mostar"
$ws.Range("G133").Value = "FK Tuzla City"
$ws.Range("H133").Value = 1
$ws.Range("I133").Value = 1
$ws.Range("J133").Value = "D"
$ws.Range("K133").Value = 1.55
$ws.Range("L133").Value = 3.6
$ws.Range("M133").Value = 6
$ws.Range("N133").Value = 1.6
$ws.Range("O133").Value = 3.5
$ws.Range("P133").Value = 5.5
$ws.Range("Q133").Value = -0.75
$ws.Range("R133").Value = 1.85
$ws.Range("S133").Value = 1.95
$ws.Range("T133").Value = 2.25
$ws.Range("U133").Value = 1.8
$ws.Range("V133").Value = 2
$ws.Range("W133").Value = -1
$ws.Range("X133").Value = 2.5
$ws.Range("Y133").Value = -1
$ws.Range("Z133").Value = -1
$ws.Range("AA133").Value = 0.95
$ws.Range("AB133").Value = -0.5
$ws.Range("AC133").Value = 0.5
$ws.Range("A133").Font.Bold = $true
$ws.Range("A133").HorizontalAlignment = -4108
$ws.Range("A133").VerticalAlignment = -4160
$ws.Range("A133").Borders.LineStyle = 1
$ws.Range("E133").NumberFormat = "YYYY-MM-DD HH:MM:SS"
